$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97; everything from the old row 97 downward
# shifts down by one (old row 97 -> new row 98, ..., old row 142 -> new row 143).
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new weekly record.
$ws.Range("A97").Value = 7
$ws.Range("B97").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C97").Value = "Ñuble"
$ws.Range("D97").Value = Get-Date -Year 2023 -Month 6 -Day 12 -Hour 0 -Minute 0 -Second 0
$ws.Range("E97").Value = 16
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100108
$ws.Range("H97").Value = "Tropicales y subtropicales"
$ws.Range("I97").Value = 100108002
$ws.Range("J97").Value = "Mango"
$ws.Range("K97").Value = "Sin especificar"
$ws.Range("L97").Value = "Primera"
$ws.Range("M97").Value = 70
$ws.Range("N97").Value = 9000
$ws.Range("O97").Value = 10000
$ws.Range("P97").Value = 9571
$ws.Range("Q97").Value = "$/bandeja 4 kilos"
$ws.Range("R97").Value = "Perú"
$ws.Range("S97").Value = 2393
$ws.Range("T97").Value = 4
